# Auto-generated Excel COM-interop script to apply the commit diff
# Updates leve-profit calculation columns (H-N) across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3458.0322
$ws.Range("I64").Value = 3344.4443
$ws.Range("J64").Value = 3615.3076
$ws.Range("K64").Value = 3344.4443
$ws.Range("L64").Value = 3615.3076
$ws.Range("M64").Value = -3096.4443
$ws.Range("N64").Value = -4111.3076
# Row 67
$ws.Range("H67").Value = 3458.0322
$ws.Range("I67").Value = 3344.4443
$ws.Range("J67").Value = 3615.3076
$ws.Range("K67").Value = 3344.4443
$ws.Range("L67").Value = 3615.3076
$ws.Range("M67").Value = -2486.4443
$ws.Range("N67").Value = -5331.3076
# Row 69
$ws.Range("H69").Value = 22247296
$ws.Range("I69").Value = 3100
$ws.Range("J69").Value = 23698004
$ws.Range("K69").Value = 9300
$ws.Range("L69").Value = 71094012
$ws.Range("M69").Value = -8426
$ws.Range("N69").Value = -71095760
# Row 72
$ws.Range("H72").Value = 22247296
$ws.Range("I72").Value = 3100
$ws.Range("J72").Value = 23698004
$ws.Range("K72").Value = 27900
$ws.Range("L72").Value = 213282036
$ws.Range("M72").Value = -23532
$ws.Range("N72").Value = -213290772
# Row 138
$ws.Range("H138").Value = 2833.3225
$ws.Range("I138").Value = 1171.0312
$ws.Range("J138").Value = 4606.433
$ws.Range("K138").Value = 3513.0936
$ws.Range("L138").Value = 13819.299
$ws.Range("M138").Value = 1626.9064
$ws.Range("N138").Value = -24099.299

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17245770
$ws.Range("I32").Value = 19233612
$ws.Range("J32").Value = 17802.166
$ws.Range("K32").Value = 19233612
$ws.Range("L32").Value = 17802.166
$ws.Range("M32").Value = -19233325
$ws.Range("N32").Value = -18376.166
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
# Row 110
$ws.Range("H110").Value = 1405.6538
$ws.Range("I110").Value = 1184.8636
$ws.Range("J110").Value = 2620
$ws.Range("K110").Value = 1184.8636
$ws.Range("L110").Value = 2620
$ws.Range("M110").Value = 860.1364000000001
$ws.Range("N110").Value = -6710

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 820
$ws.Range("I94").Value = 820
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 820
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -369
$ws.Range("N94").Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1485.9362
$ws.Range("I58").Value = 949.6857
$ws.Range("J58").Value = 3050
$ws.Range("K58").Value = 949.6857
$ws.Range("L58").Value = 3050
$ws.Range("M58").Value = -746.6857
$ws.Range("N58").Value = -3456
# Row 96
$ws.Range("H96").Value = 7058.1665
$ws.Range("J96").Value = 7058.1665
$ws.Range("L96").Value = 7058.1665
$ws.Range("N96").Value = -12550.1665
# Row 134
$ws.Range("H134").Value = 2360.35
$ws.Range("I134").Value = 1121.8334
$ws.Range("J134").Value = 13507
$ws.Range("K134").Value = 3365.5002
$ws.Range("L134").Value = 40521
$ws.Range("M134").Value = -830.5001999999999
$ws.Range("N134").Value = -45591
# Row 136
$ws.Range("H136").Value = 1485.9362
$ws.Range("I136").Value = 949.6857
$ws.Range("J136").Value = 3050
$ws.Range("K136").Value = 2849.0571
$ws.Range("L136").Value = 9150
$ws.Range("M136").Value = -299.0571
$ws.Range("N136").Value = -14250

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 584.5238000000001
$ws.Range("I5").Value = 584.5238000000001
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1753.5714
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1641.5714
$ws.Range("N5").Value = ""
# Row 99
$ws.Range("H99").Value = 2583.2222
$ws.Range("I99").Value = 1449.8
$ws.Range("K99").Value = 4349.4
$ws.Range("M99").Value = -2103.4
# Row 122
$ws.Range("H122").Value = 11528.35
$ws.Range("J122").Value = 1934
$ws.Range("L122").Value = 17406
$ws.Range("N122").Value = -22306
# Row 135
$ws.Range("H135").Value = 584.5238000000001
$ws.Range("I135").Value = 584.5238000000001
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5260.7142
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2725.7142
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 29
$ws.Range("N22").Value = -2058

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 752.58826
$ws.Range("I22").Value = 742.4167
$ws.Range("J22").Value = 777
$ws.Range("K22").Value = 742.4167
$ws.Range("L22").Value = 777
$ws.Range("M22").Value = -447.4167
$ws.Range("N22").Value = -1367
# Row 27
$ws.Range("H27").Value = 752.58826
$ws.Range("I27").Value = 742.4167
$ws.Range("J27").Value = 777
$ws.Range("K27").Value = 742.4167
$ws.Range("L27").Value = 777
$ws.Range("M27").Value = -635.4167
$ws.Range("N27").Value = -991
# Row 35
$ws.Range("H35").Value = 15744.444
$ws.Range("I35").Value = 540
$ws.Range("J35").Value = 34750
$ws.Range("K35").Value = 540
$ws.Range("L35").Value = 34750
$ws.Range("M35").Value = -204
$ws.Range("N35").Value = -35422
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = ""
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = ""
# Row 87
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32246
# Row 88
$ws.Range("H88").Value = 25085.5
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20856
# Row 90
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -101232
# Row 91
$ws.Range("H91").Value = 25085.5
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22964
# Row 111
$ws.Range("H111").Value = 47900
$ws.Range("J111").Value = 47900
$ws.Range("L111").Value = 47900
$ws.Range("N111").Value = -56080
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
# Row 128
$ws.Range("H128").Value = 59900
$ws.Range("J128").Value = 59900
$ws.Range("L128").Value = 59900
$ws.Range("N128").Value = -69860

$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Range("H114").Value = 25475
$ws.Range("J114").Value = 25475
$ws.Range("L114").Value = 25475
$ws.Range("N114").Value = -34153
